$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.263.78"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.666.50"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D5").Value = "'218.77"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'0.5230"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "'0.2664"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "'20.96"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "1.676.92"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "'4.447"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "1.892.69"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Value = "'0.5471"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "0.0₅8261"
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "'64.79"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "26.286.57"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "'4.674"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("D21").Value = "'194.28"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'10.14"
$ws.Range("E22").Value = "  -1.56%  "
$ws.Range("D23").Value = "'6.065"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'140.50"
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").Value = "'0.1235"
$ws.Range("E26").Value = "  -4.18%  "
$ws.Range("D27").Value = "'7.192"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").Value = "'1.419"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "'0.06165"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").Value = "'1.282"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "'3.588"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'3.290"
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "'0.9706"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("D36").Value = "'2.425"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "'2.790"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").Value = "'0.5727"
$ws.Range("E38").Value = "  -7.36%  "
$ws.Range("D39").Value = "'0.01611"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").Value = "'6.007"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("D41").Value = "'0.8560"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D43").Value = "1.018.71"
$ws.Range("E43").Value = "  -6.24%  "
$ws.Range("D44").Value = "'100.29"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "1.808.42"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "'57.43"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").Value = "'8.051"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "'1.490"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "'0.05188"
$ws.Range("E51").Value = "  -0.38%  "
